$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.639.59"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "1.617.00"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.990"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.10"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.06"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "1.845.78"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").Value = "1.629.95"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.556"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.60"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "27.655.98"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.23"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.64"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.04"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.39"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.88"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.42"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.16"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0478"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "1.389.27"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.840"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.50"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.82"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.35"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").Value = "1.755.54"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("E47").Value = "  -3.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.60"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.51"
$ws.Range("E51").Value = "  +0.81%  "
